$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (D) column is stored as text in the sheet even when it looks numeric
# (e.g. "27.551.81" isn't a valid number anyway, but "324.22" is and would
# otherwise be auto-converted to a real number by the Value setter). Force
# text format on the whole Price column before writing so every new value
# round-trips as inline/shared text like the original file.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 46 and 47 swap places (Decentraland <-> PancakeSwap) plus value updates.
$ws.Range("B46").Value = "PancakeSwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D46").Value = "3.711"
$ws.Range("E46").Value = "  -1.01%  "

$ws.Range("B47").Value = "Decentraland"
$ws.Range("C47").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D47").Value = "0.5870"
$ws.Range("E47").Value = "  -4.04%  "

# Per-row Price (D) and Volume(1h) (E) updates for the rest of the table.
$updates = @{
    2  = @{ D = "27.543.21";   E = "  -2.73%  " }
    3  = @{ D = "1.752.25";    E = "  -3.36%  " }
    4  = @{ E = "  +0.24%  " }
    5  = @{ D = "324.01";      E = "  -0.93%  " }
    6  = @{ E = "  +0.24%  " }
    7  = @{ D = "0.4457";      E = "  +2.08%  " }
    8  = @{ D = "0.3597";      E = "  -2.01%  " }
    9  = @{ D = "0.07500";     E = "  -2.36%  " }
    10 = @{ D = "42.22";       E = "  -5.60%  " }
    11 = @{ D = "1.101";       E = "  -3.64%  " }
    12 = @{ D = "1.003";       E = "  +0.28%  " }
    13 = @{ D = "20.69";       E = "  -6.07%  " }
    14 = @{ D = "6.034";       E = "  -4.48%  " }
    15 = @{ D = "7.198";       E = "  -4.35%  " }
    16 = @{ D = "1.758.71";    E = "  -3.82%  " }
    17 = @{ D = "92.75";       E = "  -2.85%  " }
    18 = @{ D = "0.00001063";  E = "  -1.66%  " }
    19 = @{ D = "0.06412";     E = "  -1.68%  " }
    20 = @{ D = "1.002";       E = "  +0.25%  " }
    21 = @{ D = "17.05";       E = "  -2.14%  " }
    22 = @{ D = "5.854";       E = "  -6.25%  " }
    23 = @{ D = "27.596.13";   E = "  -2.57%  " }
    24 = @{ D = "11.22";       E = "  -3.16%  " }
    25 = @{ D = "2.098";       E = "  +0.49%  " }
    26 = @{ D = "162.51";      E = "  +0.32%  " }
    27 = @{ D = "20.43";       E = "  -1.35%  " }
    28 = @{ D = "1.956.94";    E = "  -3.44%  " }
    29 = @{ D = "2.125";       E = "  -6.74%  " }
    30 = @{ D = "125.88";      E = "  -2.40%  " }
    31 = @{ D = "1.090";       E = "  -9.66%  " }
    32 = @{ D = "0.09045";     E = "  -1.54%  " }
    33 = @{ D = "3.639";       E = "  +3.66%  " }
    34 = @{ D = "5.544";       E = "  -7.22%  " }
    35 = @{ D = "12.04";       E = "  -7.52%  " }
    36 = @{ D = "0.02297";     E = "  -2.15%  " }
    37 = @{ D = "0.2102";      E = "  -3.21%  " }
    38 = @{ D = "0.6386";      E = "  -2.99%  " }
    39 = @{ D = "0.05975";     E = "  -3.81%  " }
    40 = @{ D = "4.927";       E = "  -5.17%  " }
    41 = @{ D = "1.193";       E = "  -0.23%  " }
    42 = @{ E = "  +0.28%  " }
    43 = @{ D = "1.392";       E = "  -2.50%  " }
    44 = @{ D = "7.783";       E = "  -4.22%  " }
    45 = @{ D = "13.20";       E = "  -5.51%  " }
    48 = @{ D = "1.959";       E = "  -2.85%  " }
    49 = @{ D = "121.61";      E = "  -3.41%  " }
    50 = @{ D = "1.155";       E = "  -0.10%  " }
    51 = @{ D = "0.06872";     E = "  -1.85%  " }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
